# Contacts - 16 Dec 2024
# - Users!A2: rename contact "Ayati Arvind" -> "Amanda Donovan"
# - AffiliatedCompany sheet: move the saved selection from C1:E2 to E20
# - Users sheet: move the saved selection from D8 to A2 (and keep it the
#   active/visible tab when the workbook is saved)

$wb = $excel.ActiveWorkbook

# Touch the AffiliatedCompany sheet's selection first so it does not end up
# being the last-activated (and therefore "tabSelected") sheet.
$wsAffiliated = $wb.Worksheets.Item("AffiliatedCompany")
$wsAffiliated.Activate()
$wsAffiliated.Range("E20").Select()

# Now switch to Users, update the contact name, and leave the selection/tab
# focus here, matching the saved workbook state.
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Activate()
$wsUsers.Range("A2").Value = "Amanda Donovan"
$wsUsers.Range("A2").Select()
